# Update the "partidas_mines" data table in Hoja1.
# New data shifts "partida" numbers from 30-40 up to 40-47 and reduces the
# number of data rows from 44 (rows 2-45) down to 32 (rows 2-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @(40,1,1),
    @(40,3,2),
    @(40,3,5),
    @(40,4,1),
    @(41,1,1),
    @(41,2,4),
    @(41,4,1),
    @(41,4,4),
    @(42,1,2),
    @(42,1,4),
    @(42,3,1),
    @(42,5,3),
    @(43,1,1),
    @(43,2,2),
    @(43,3,1),
    @(43,3,4),
    @(44,2,3),
    @(44,2,4),
    @(44,4,5),
    @(44,5,2),
    @(45,1,1),
    @(45,1,3),
    @(45,1,4),
    @(45,4,1),
    @(46,1,3),
    @(46,3,3),
    @(46,4,4),
    @(46,4,5),
    @(47,1,2),
    @(47,4,1),
    @(47,4,3),
    @(47,5,4)
)

# Clear out the previous data range (rows 2-45) before writing the new,
# shorter table (rows 2-33).
$ws.Range("A2:C45").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Column F picked up a wide "best fit" width in the saved workbook.
$ws.Columns.Item(6).ColumnWidth = 88.7109375

# Match the recorded selection from the edited file.
$ws.Range("F26").Select()
